# Updating daily files July 20
# The daily refresh re-sorted a handful of pitcher blocks (Kyle Gibson, Charlie
# Morton, Taj Bradley, Nestor Cortes). Each pitcher occupies two rows (vs L /
# vs R splits). The net effect is that the 4-row block for Taj Bradley/Nestor
# Cortes (previously rows 308-311) now swaps places with the 4-row block for
# Kyle Gibson/Charlie Morton (previously rows 304-307), while everything else
# on the sheet stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rngA = $ws.Range("A304:AA307")
$rngB = $ws.Range("A308:AA311")

$valsA = $rngA.Value2
$valsB = $rngB.Value2

$rngA.Value2 = $valsB
$rngB.Value2 = $valsA
